$wb = $excel.ActiveWorkbook

# The new "UK" sheet mirrors the existing "Czech" market sheet (same
# layout / styles / merged cells / repeater list), so clone that sheet
# and drop the copy at the very end of the tab strip (after "Poland").
$src = $wb.Worksheets.Item("Czech")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$src.Copy($null, $lastSheet)

$ws = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Name = "UK"

# Match the column widths used by the other "standard" market sheets.
$ws.Columns.Item(1).ColumnWidth = 24.77734375
$ws.Columns.Item(2).ColumnWidth = 15.21875
$ws.Columns.Item(3).ColumnWidth = 13.44140625
$ws.Columns.Item(4).ColumnWidth = 16.44140625

# Fill in the market-specific values (B4 is entered before B2, matching
# the shared-string insertion order of the source edit).
$ws.Range("B4").Value = "NGC-2741/T3343"
$ws.Range("B2").Value = "UKl Market"

# Leave the selection on B4, as in the authored sheet.
[void]$ws.Range("B4").Select()
